$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed assignments first, in the exact order needed so that the
# shared-string table ends up with entries 47-53 in the same order
# as the target workbook (Excel assigns new shared-string ids in
# first-use order).
$ws.Range("E2").Value = "Network Added to Base Service"
$ws.Range("E6").Value = "Network Removed from Add-On Package"
$ws.Range("E7").Value = "Network Removed from Database"
$ws.Range("E15").Value = "Name of Add-On Package Changed"
$ws.Range("E24").Value = "Network Added to Add-On Package"
$ws.Range("E9").Value = "New Network Added to Database in Jun 2020"
$ws.Range("E14").Value = "Network Removed from Base Service"

# Remaining rows (values already exist in the shared string table)
$ws.Range("E3").Value = "Network Added to Base Service"
$ws.Range("E4").Value = "Network Added to Base Service"
$ws.Range("E5").Value = "Network Added to Base Service"
$ws.Range("E8").Value = "Network Removed from Database"
$ws.Range("E10").Value = "New Network Added to Database in Jun 2020"
$ws.Range("E11").Value = "New Network Added to Database in Jun 2020"
$ws.Range("E12").Value = "New Network Added to Database in Jun 2020"
$ws.Range("E13").Value = "New Network Added to Database in Jun 2020"
$ws.Range("E16").Value = "Name of Add-On Package Changed"
$ws.Range("E17").Value = "Name of Add-On Package Changed"
$ws.Range("E18").Value = "Network Removed from Base Service"
$ws.Range("E19").Value = "Network Removed from Add-On Package"
$ws.Range("E20").Value = "Name of Add-On Package Changed"
$ws.Range("E21").Value = "Name of Add-On Package Changed"
$ws.Range("E22").Value = "Name of Add-On Package Changed"
$ws.Range("E23").Value = "Network Removed from Base Service"
$ws.Range("E25").Value = "Network Added to Add-On Package"
$ws.Range("E26").Value = "Network Added to Add-On Package"
$ws.Range("E27").Value = "Network Added to Add-On Package"
$ws.Range("E28").Value = "Network Added to Add-On Package"
$ws.Range("E29").Value = "Network Added to Add-On Package"
$ws.Range("E30").Value = "Network Removed from Base Service"
$ws.Range("E31").Value = "Network Removed from Base Service"
$ws.Range("E32").Value = "Network Removed from Add-On Package"
$ws.Range("E33").Value = "Network Removed from Base Service"
$ws.Range("E34").Value = "Network Added to Add-On Package"
$ws.Range("E35").Value = "Network Added to Add-On Package"

$wb.Save()
